# Update the "Förändrad" date column (C) for rows 2-6 from 45208 to 45212
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
